$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 22.71826066666667
$ws.Range("H2").Value = 68.154782
$ws.Range("I2").Value = 0.871041851052479
$ws.Range("J2").Value = 0.871041851052479
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 93.85711466666667
$ws.Range("N2").Value = 281.571344
$ws.Range("Q2").Value = 2132.270396418557
$ws.Range("R2").Value = 19190.43356776701
$ws.Range("S2").Value = 0.871041851052479
$ws.Range("T2").Value = 0.871041851052479

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.436847666666667
$ws.Range("H3").Value = 7.310543
$ws.Range("I3").Value = 0.09343128567146972
$ws.Range("J3").Value = 0.09343128567146973
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 93.85711466666667
$ws.Range("N3").Value = 281.571344
$ws.Range("Q3").Value = 228.7154908755325
$ws.Range("R3").Value = 2058.439417879792
$ws.Range("S3").Value = 0.09343128567146972
$ws.Range("T3").Value = 0.09343128567146973

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.9266013333333333
$ws.Range("H4").Value = 2.779804
$ws.Range("I4").Value = 0.03552686327605135
$ws.Range("J4").Value = 0.03552686327605135
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 93.85711466666667
$ws.Range("N4").Value = 281.571344
$ws.Range("Q4").Value = 86.9681275929529
$ws.Range("R4").Value = 782.7131483365761
$ws.Range("S4").Value = 0.03552686327605135
$ws.Range("T4").Value = 0.03552686327605135
